# "buck changed and calculated"
#
# Extend the "buck" sheet buck-converter calculator with extra parameters
# (V_REF / R_FBB feedback-resistor pair, inductor value L, enable resistors
# R_ENB / R_ENT) and new derived formulas (the L_MIN echo in K11, the R_FBT
# feedback-resistor formula in G18 and the ripple-current formula in K18),
# boxing all the input/output cells with a thin border.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("buck")

# ------------------------------------------------------------------
# 1) New label / unit strings, entered in the exact order the sheet
#    first needs them so the shared-string table comes out in the
#    same sequence as in the authored workbook.
# ------------------------------------------------------------------
$ws.Range("F11").Value = "V_REF"
$ws.Range("F12").Value = "R_FBB"
$ws.Range("L11").Value = "[H]"
$ws.Range("H11").Value = "[V]"
$ws.Range("P11").Value = "[Ohm]"
$ws.Range("D11").Value = "[A]"
$ws.Range("D14").Value = "[Hz]"
$ws.Range("J18").Value = "∆ iL"
$ws.Range("J11").Value = "L"
$ws.Range("N11").Value = "R_ENB"
$ws.Range("N12").Value = "R_ENT"

# Cells that reuse one of the strings just created above.
$ws.Range("D12").Value = "[V]"
$ws.Range("D13").Value = "[V]"
$ws.Range("H12").Value = "[Ohm]"
$ws.Range("P12").Value = "[Ohm]"
$ws.Range("H18").Value = "[Ohm]"
$ws.Range("D18").Value = "[H]"

# ------------------------------------------------------------------
# 2) New numeric inputs / blank placeholder cells.
# ------------------------------------------------------------------
$ws.Range("G11").Value = 0.8
$ws.Range("G12").Value = 22100
$ws.Range("O11").Value = ""
$ws.Range("O12").Value = ""
$ws.Range("D15").Value = ""

# ------------------------------------------------------------------
# 3) New formulas.
# ------------------------------------------------------------------
$ws.Range("K11").Formula = "=C18"
$ws.Range("G18").Formula = "=((C12-G11)/G11)*G12"
$ws.Range("K18").Formula = "=(C12*(C13-C12))/(C13*K11*C14)"

# ------------------------------------------------------------------
# 4) Remove the now-obsolete helper cell C21 (its value now lives in
#    K11), leaving a couple of blank formatted placeholder cells
#    behind (C20, C22) as in the new layout.
# ------------------------------------------------------------------
$ws.Range("C21").Clear() | Out-Null
$ws.Range("C20").Value = ""
$ws.Range("C22").Value = ""

# ------------------------------------------------------------------
# 5) Number formats.
# ------------------------------------------------------------------
$ws.Range("K11").NumberFormat = "_-* #,##0.0000000000_-;\-* #,##0.0000000000_-;_-* ""-""??_-;_-@_-"
$ws.Range("K18").NumberFormat = "_-* #,##0.00\ _z_ł_-;\-* #,##0.00\ _z_ł_-;_-* ""-""??????????\ _z_ł_-;_-@_-"
$ws.Range("C20").NumberFormat = "0.00E+00"
$ws.Range("C22").NumberFormat = "0.00E+00"

# ------------------------------------------------------------------
# 6) Fonts: every label cell in the new grid uses the bold label font
#    that B11:B18 already used; apply it to the new F/J/N label cells
#    too.
# ------------------------------------------------------------------
$boldLabels = "F11,F12,F18,J11,J18,N11,N12"
foreach ($addr in $boldLabels.Split(",")) {
    $ws.Range($addr).Font.Bold = $true
}

# ------------------------------------------------------------------
# 7) Alignment: unit-label cells are centred, C15 stays right aligned.
# ------------------------------------------------------------------
$unitCells = "D11,D12,D13,D14,D18,H11,H12,H18,L11,P11,P12"
foreach ($addr in $unitCells.Split(",")) {
    $ws.Range($addr).HorizontalAlignment = -4108
    $ws.Range($addr).VerticalAlignment = -4108
}
$ws.Range("C15").HorizontalAlignment = -4152

# ------------------------------------------------------------------
# 8) Thin box border around every cell of the new input/output grid.
# ------------------------------------------------------------------
$boxedRanges = "B11:D11,F11:H11,J11:L11,N11:P11,B12:D12,F12:H12,N12:P12,B13:D13,B14:D14,B15:D15,B18:D18,F18:H18,J18:K18"
foreach ($addr in $boxedRanges.Split(",")) {
    $rng = $ws.Range($addr)
    $rng.Borders.LineStyle = 1
    $rng.Borders.Weight = 2
}

# ------------------------------------------------------------------
# 9) Sheet cosmetics to mirror the updated view/dimension.
# ------------------------------------------------------------------
$ws.Columns.Item(11).ColumnWidth = 13.3
$ws.Range("M28").Select() | Out-Null

$wb.Save()
